# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Vega Monumental
# Concepción / Frutilla dataset, pushing the existing records down by
# two rows (row 132 -> 134, etc.), and populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 132, shifting rows 132:246 down to 134:248.
$ws.Rows.Item(132).Insert()
$ws.Rows.Item(132).Insert()

# New row 132
$ws.Cells.Item(132, 1).Value = 11
$ws.Cells.Item(132, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(132, 3).Value = "Bíobío"
$ws.Cells.Item(132, 4).Value = 44566
$ws.Cells.Item(132, 5).Value = 8
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100101
$ws.Cells.Item(132, 8).Value = "Berries"
$ws.Cells.Item(132, 9).Value = 100112025
$ws.Cells.Item(132, 10).Value = "Frutilla"
$ws.Cells.Item(132, 11).Value = "Sin especificar"
$ws.Cells.Item(132, 12).Value = "Especial"
$ws.Cells.Item(132, 13).Value = 200
$ws.Cells.Item(132, 14).Value = 7000
$ws.Cells.Item(132, 15).Value = 7500
$ws.Cells.Item(132, 16).Value = 7250
$ws.Cells.Item(132, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(132, 18).Value = "Región del Maule"
$ws.Cells.Item(132, 19).Value = 1036
$ws.Cells.Item(132, 20).Value = 7

# New row 133
$ws.Cells.Item(133, 1).Value = 11
$ws.Cells.Item(133, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(133, 3).Value = "Bíobío"
$ws.Cells.Item(133, 4).Value = 44566
$ws.Cells.Item(133, 5).Value = 8
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100101
$ws.Cells.Item(133, 8).Value = "Berries"
$ws.Cells.Item(133, 9).Value = 100112025
$ws.Cells.Item(133, 10).Value = "Frutilla"
$ws.Cells.Item(133, 11).Value = "Sin especificar"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 250
$ws.Cells.Item(133, 14).Value = 6000
$ws.Cells.Item(133, 15).Value = 6500
$ws.Cells.Item(133, 16).Value = 6260
$ws.Cells.Item(133, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(133, 18).Value = "Región del Maule"
$ws.Cells.Item(133, 19).Value = 894
$ws.Cells.Item(133, 20).Value = 7
